$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

for ($r = 2; $r -le 379; $r++) {
    $cell = $ws.Cells.Item($r, 6)
    $v = $cell.Value2
    if ($v -eq 1) {
        $cell.Value = "1T"
    } elseif ($v -eq 2) {
        $cell.Value = "2T"
    }
}

$ws.Range("F6:F378").Select()
